$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ E=3; G=468.2935736666666; H=1404.880721;       I=0.2378793178672371; J=0.2378793178672371; K=3; M=261.380203;        N=784.1406089999999;  O=0.6968677182772199;  P=0.6968677182772199;  Q=122402.6693485888; R=1101624.024137299; S=0.1657704174674831;  T=0.165770417467483 }
    3  = @{ E=3; G=468.2935736666666; H=1404.880721;       I=0.2378793178672371; J=0.2378793178672371; K=3; M=31.999428;          N=95.998284;           O=0.08531391482826334; P=0.08531391482826335; Q=14985.12649340919; R=134866.1384406828; S=0.02029441586393085; T=0.02029441586393085 }
    4  = @{ E=3; G=468.2935736666666; H=1404.880721;       I=0.2378793178672371; J=0.2378793178672371; K=3; M=81.69901900000001;  N=245.097057;          O=0.2178183668945166;  P=0.2178183668945167;  Q=38259.1255725709;  R=344332.1301531381; S=0.05181448453582321; T=0.05181448453582321 }
    5  = @{ E=3; G=715.8492226666667; H=2147.547668;        I=0.3636302831371944; J=0.3636302831371944; K=3; M=261.380203;        N=784.1406089999999;  O=0.6968677182772199;  P=0.6968677182772199;  Q=187108.8151380055; R=1683979.33624205;  S=0.2534022057063161;  T=0.2534022057063161 }
    6  = @{ E=3; G=715.8492226666667; H=2147.547668;        I=0.3636302831371944; J=0.3636302831371944; K=3; M=31.999428;          N=95.998284;           O=0.08531391482826334; P=0.08531391482826335; Q=22906.76565957797; R=206160.8909362017; S=0.03102272300454389; T=0.03102272300454389 }
    7  = @{ E=3; G=715.8492226666667; H=2147.547668;        I=0.3636302831371944; J=0.3636302831371944; K=3; M=81.69901900000001;  N=245.097057;          O=0.2178183668945166;  P=0.2178183668945167;  Q=58484.17924377924; R=526357.6131940131; S=0.07920535442633438; T=0.07920535442633438 }
    8  = @{ E=3; G=784.4754839999999; H=2353.426452;        I=0.3984903989955685; J=0.3984903989955685; K=3; M=261.380203;        N=784.1406089999999;  O=0.6968677182772199;  P=0.6968677182772199;  Q=205046.3612564432; R=1845417.251307989; S=0.2776950951034208;  T=0.2776950951034208 }
    9  = @{ E=3; G=784.4754839999999; H=2353.426452;        I=0.3984903989955685; J=0.3984903989955685; K=3; M=31.999428;          N=95.998284;           O=0.08531391482826334; P=0.08531391482826335; Q=25102.76676802315; R=225924.9009122083; S=0.03399677595978861; T=0.03399677595978861 }
    10 = @{ E=3; G=784.4754839999999; H=2353.426452;        I=0.3984903989955685; J=0.3984903989955685; K=3; M=81.69901900000001;  N=245.097057;          O=0.2178183668945166;  P=0.2178183668945167;  Q=64090.8774723502;  R=576817.8972511517; S=0.08679852793235908; T=0.08679852793235908 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
